$wb = $excel.ActiveWorkbook

# --- DeviceInfo sheet: new device under test has a different IP / calibration file
$wsDevice = $wb.Worksheets.Item("DeviceInfo")
$wsDevice.Activate()
$wsDevice.Range("A1:XFD1048576").Select()
$wsDevice.Range("A2").Value = "10.75.58.66"
$wsDevice.Range("B2").Value = "6U_10I.cal"

# --- Cabling sheet: drop the two extra placeholder channel rows (17 & 18) that
#     referenced the now-removed Channel[16]/Channel[17] labels. Keep the W
#     (lookup list) column intact, only clear the data columns A:J.
$wsCabling = $wb.Worksheets.Item("Cabling")
$wsCabling.Activate()
$wsCabling.Range("A1:XFD1048576").Select()
$wsCabling.Range("A18:J19").ClearContents()

# Shrink the data-validation ranges so they no longer cover the cleared rows.
# Deleting validation on just the trailing rows causes Excel to shrink the
# sqref of the remaining rule rather than creating a brand-new one.
$wsCabling.Range("C18:C19").Validation.Delete()
$wsCabling.Range("E18:E19").Validation.Delete()
$wsCabling.Range("F18:F19").Validation.Delete()
$wsCabling.Range("G18:G19").Validation.Delete()
$wsCabling.Range("J18:J19").Validation.Delete()
$wsCabling.Range("B19:B20").Validation.Delete()

# --- BusbarFeederMap sheet: updated feeder/channel counts for the new cabling
$wsBusbar = $wb.Worksheets.Item("BusbarFeederMap")
$wsBusbar.Activate()
$wsBusbar.Range("A1:XFD1048576").Select()
$wsBusbar.Range("B37").Value = 0
$wsBusbar.Range("B38").Value = 0
$wsBusbar.Range("B40").Value = 1

# --- DSPChannelMap sheet: updated DSP2 feeder map values
$wsDsp = $wb.Worksheets.Item("DSPChannelMap")
$wsDsp.Activate()
$wsDsp.Range("A1:XFD1048576").Select()
$wsDsp.Range("D9").Value = 0
$wsDsp.Range("D10").Value = 0

# --- DSPFeederMap sheet: touch selection so the full-sheet selection is saved
#     consistently across every sheet (matches every other tab in the file).
$wsFeeder = $wb.Worksheets.Item("DSPFeederMap")
$wsFeeder.Activate()
$wsFeeder.Range("A1:XFD1048576").Select()

# Leave DSPChannelMap as the active tab, matching the saved file.
$wsDsp.Activate()
$wsDsp.Range("A1:XFD1048576").Select()
